$wb = $excel.ActiveWorkbook

function Set-Cells {
    param($ws, $row, $values)
    foreach ($col in $values.Keys) {
        $ws.Range("$col$row").Value = $values[$col]
    }
}

$ws = $wb.Worksheets.Item("ALC")
Set-Cells -ws $ws -row 17 -values @{ "H" = 877.4286; "I" = 381.76923; "J" = 1307; "K" = 1145.30769; "L" = 3921; "M" = -977.3076900000001; "N" = -4257 }
Set-Cells -ws $ws -row 19 -values @{ "H" = 26187; "I" = 2000; "J" = 29642.285; "K" = 2000; "L" = 29642.285; "M" = -1825; "N" = -29992.285 }
Set-Cells -ws $ws -row 33 -values @{ "H" = 96.933334; "I" = 102.23077; "J" = 62.5; "K" = 102.23077; "L" = 62.5; "M" = 126.76923; "N" = -520.5 }
Set-Cells -ws $ws -row 40 -values @{ "H" = 1333.3334; "I" = 1500; "J" = 1250; "K" = 1500; "L" = 1250; "M" = -1325; "N" = -1600 }
Set-Cells -ws $ws -row 42 -values @{ "H" = 596.5; "I" = 90; "J" = 849.75; "K" = 270; "L" = 2549.25; "M" = -40; "N" = -3009.25 }
Set-Cells -ws $ws -row 127 -values @{ "H" = 917.5263; "I" = 361; "J" = 1682.75; "K" = 1083; "L" = 5048.25; "M" = 3877; "N" = -14968.25 }
Set-Cells -ws $ws -row 137 -values @{ "H" = 1513.2368; "I" = 1334.5416; "J" = 1819.5714; "K" = 4003.6248; "L" = 5458.7142; "M" = -1453.6248; "N" = -10558.7142 }
$ws = $wb.Worksheets.Item("ARM")
Set-Cells -ws $ws -row 23 -values @{ "H" = 15222.556; "I" = 70006; "J" = 8374.625; "K" = 70006; "L" = 8374.625; "M" = -69747; "N" = -8892.625 }
Set-Cells -ws $ws -row 61 -values @{ "H" = 3481.5; "I" = 3493; "J" = 3473.0667; "K" = 3493; "L" = 3473.0667; "M" = -3281; "N" = -3897.0667 }
Set-Cells -ws $ws -row 74 -values @{ "H" = 1244.258; "I" = 967.4706; "J" = 1580.3572; "K" = 967.4706; "L" = 1580.3572; "M" = -93.47059999999999; "N" = -3328.3572 }
Set-Cells -ws $ws -row 77 -values @{ "H" = 1244.258; "I" = 967.4706; "J" = 1580.3572; "K" = 4837.353; "L" = 7901.786; "M" = -469.3530000000001; "N" = -16637.786 }
Set-Cells -ws $ws -row 122 -values @{ "H" = 44479.824; "I" = 56406; "J" = 1545.6; "K" = 169218; "L" = 4636.799999999999; "M" = -166768; "N" = -9536.799999999999 }
Set-Cells -ws $ws -row 132 -values @{ "H" = 3025.322; "I" = 2455.5854; "J" = 4323.0557; "K" = 7366.7562; "L" = 12969.1671; "M" = -4836.7562; "N" = -18029.1671 }
Set-Cells -ws $ws -row 136 -values @{ "H" = 3481.5; "I" = 3493; "J" = 3473.0667; "K" = 10479; "L" = 10419.2001; "M" = -7929; "N" = -15519.2001 }
$ws = $wb.Worksheets.Item("BSM")
Set-Cells -ws $ws -row 22 -values @{ "H" = 274.2857; "I" = 274.2857; "J" = 0; "K" = 274.2857; "L" = 0; "M" = -101.2857 }
Set-Cells -ws $ws -row 134 -values @{ "H" = 3224.6765; "I" = 2791.6667; "J" = 3924.1538; "K" = 8375.000100000001; "L" = 11772.4614; "M" = -5840.000100000001; "N" = -16842.4614 }
Set-Cells -ws $ws -row 137 -values @{ "H" = 50780; "I" = 0; "J" = 50780; "K" = 0; "L" = 50780; "N" = -60980 }
$ws = $wb.Worksheets.Item("CRP")
Set-Cells -ws $ws -row 31 -values @{ "H" = 4439.696; "I" = 1224.8334; "J" = 7946.8184; "K" = 1224.8334; "L" = 7946.8184; "M" = -929.8334; "N" = -8536.8184 }
Set-Cells -ws $ws -row 34 -values @{ "H" = 4439.696; "I" = 1224.8334; "J" = 7946.8184; "K" = 1224.8334; "L" = 7946.8184; "M" = -1022.8334; "N" = -8350.8184 }
Set-Cells -ws $ws -row 58 -values @{ "H" = 1218.475; "I" = 965.26086; "J" = 1561.0588; "K" = 965.26086; "L" = 1561.0588; "M" = -762.26086; "N" = -1967.0588 }
Set-Cells -ws $ws -row 118 -values @{ "H" = 38333; "I" = 0; "J" = 38333; "K" = 0; "L" = 38333; "N" = -41647 }
Set-Cells -ws $ws -row 136 -values @{ "H" = 1218.475; "I" = 965.26086; "J" = 1561.0588; "K" = 2895.78258; "L" = 4683.1764; "M" = -345.7825800000001; "N" = -9783.1764 }
$ws = $wb.Worksheets.Item("CUL")
Set-Cells -ws $ws -row 32 -values @{ "H" = 10419792; "I" = 0; "J" = 10419792; "K" = 0; "L" = 31259376; "N" = -31259942 }
Set-Cells -ws $ws -row 70 -values @{ "H" = 1368.6666; "I" = 842.4; "J" = 4000; "K" = 2527.2; "L" = 12000; "M" = -2212.2; "N" = -12630 }
Set-Cells -ws $ws -row 73 -values @{ "H" = 1368.6666; "I" = 842.4; "J" = 4000; "K" = 2527.2; "L" = 12000; "M" = -1435.2; "N" = -14184 }
Set-Cells -ws $ws -row 80 -values @{ "H" = 2988; "I" = 2501; "J" = 3312.6667; "K" = 7503; "L" = 9938.000100000001; "M" = -6567; "N" = -11810.0001 }
Set-Cells -ws $ws -row 83 -values @{ "H" = 2988; "I" = 2501; "J" = 3312.6667; "K" = 22509; "L" = 29814.0003; "M" = -17829; "N" = -39174.0003 }
$ws = $wb.Worksheets.Item("GSM")
Set-Cells -ws $ws -row 95 -values @{ "H" = 75922; "I" = 0; "J" = 75922; "K" = 0; "L" = 75922; "N" = -81414 }
Set-Cells -ws $ws -row 126 -values @{ "H" = 1657.75; "I" = 1558.5; "J" = 1757; "K" = 4675.5; "L" = 5271; "M" = -2205.5; "N" = -10211 }
Set-Cells -ws $ws -row 132 -values @{ "H" = 2064.9473; "I" = 1828.2727; "J" = 2390.375; "K" = 5484.8181; "L" = 7171.125; "M" = -2954.8181; "N" = -12231.125 }
$ws = $wb.Worksheets.Item("LTW")
Set-Cells -ws $ws -row 132 -values @{ "H" = 2095.8442; "I" = 2149.76; "J" = 1996; "K" = 6449.280000000001; "L" = 5988; "M" = -3919.280000000001; "N" = -11048 }
Set-Cells -ws $ws -row 136 -values @{ "H" = 7577464.5; "I" = 1813.8667; "J" = 23811000; "K" = 5441.6001; "L" = 71433000; "M" = -2891.6001; "N" = -71438100 }
$ws = $wb.Worksheets.Item("WVR")
Set-Cells -ws $ws -row 14 -values @{ "H" = 8840801; "I" = 67000; "J" = 14690002; "K" = 67000; "L" = 14690002; "M" = -66832; "N" = -14690338 }
Set-Cells -ws $ws -row 70 -values @{ "H" = 42975; "I" = 0; "J" = 42975; "K" = 0; "L" = 42975; "N" = -43605 }
Set-Cells -ws $ws -row 73 -values @{ "H" = 42975; "I" = 0; "J" = 42975; "K" = 0; "L" = 42975; "N" = -45159 }
Set-Cells -ws $ws -row 92 -values @{ "H" = 90550; "I" = 0; "J" = 90550; "K" = 0; "L" = 90550; "N" = -95542 }
Set-Cells -ws $ws -row 132 -values @{ "H" = 2874547; "I" = 898.55; "J" = 9260432; "K" = 2695.65; "L" = 27781296; "M" = -165.6499999999996; "N" = -27786356 }
Set-Cells -ws $ws -row 136 -values @{ "H" = 2391.7163; "I" = 2172.3; "J" = 3037.0588; "K" = 6516.900000000001; "L" = 9111.1764; "M" = -3966.900000000001; "N" = -14211.1764 }
